# "Generate Report for Handoff": refresh the "Latest Handoff Date/Datetime"
# timestamp for every localization row that is currently awaiting handoff
# (status "Handback transform failed" or "Ready for handoff"). Rows that are
# still "In Translation" or already "Handed back" keep their existing dates.
#
# New timestamps, per-sheet:
#   Overview ("Latest Handoff Date", col D): 2016-03-20 17:48:28
#   zh-cn    ("Latest Handoff Datetime", col E): 2016-03-20 17:48:19
#   de-de    ("Latest Handoff Datetime", col E): 2016-03-20 17:48:28

$wb = $excel.ActiveWorkbook

$rowsToRefresh = @(4, 6, 7, 8, 9, 10)

$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rowsToRefresh) {
    $wsOverview.Cells.Item($r, 4).Value = "2016-03-20 17:48:28"
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rowsToRefresh) {
    $wsZhCn.Cells.Item($r, 5).Value = "2016-03-20 17:48:19"
}

$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rowsToRefresh) {
    $wsDeDe.Cells.Item($r, 5).Value = "2016-03-20 17:48:28"
}
